# ContactPage.xlsx update — add PostalCode/State/PrivacyConsent locator rows,
# switch calculation back to automatic, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Workbook was saved with manual calculation; restore automatic calculation
# (this drops the calcMode="manual" attribute from <calcPr>).
$excel.Calculation = -4105   # xlCalculationAutomatic

# New locator rows appended after the existing 21 data rows (row 22 was the
# last one, new rows are 23-28).
$newRows = @(
    @{ Row = 23; Sno = 22; Name = "ContactPage_TextBox_Zip";                                         Value = "input#PostalCode";                    Styled = $false },
    @{ Row = 24; Sno = 23; Name = "ContactPage_ErrorMessage_InvalidZip_TextBox_Zip";                  Value = ".mktoError #ValidMsgPostalCode";       Styled = $true  },
    @{ Row = 25; Sno = 24; Name = "ContactPage_Dropdown_Province";                                    Value = "select#State";                        Styled = $false },
    @{ Row = 26; Sno = 25; Name = "ContactPage_ErrorMessage_InvalidProvince_Dropdown_Province";       Value = ".mktoError #ValidMsgState";            Styled = $true  },
    @{ Row = 27; Sno = 26; Name = "ContactPage_CheckBox_PrivacyConsent";                              Value = "input#mktoCheckbox_142098_0";          Styled = $true  },
    @{ Row = 28; Sno = 27; Name = "ContactPage_Label_PrivacyConsent";                                 Value = "label#LblmktoCheckbox_142098_0";       Styled = $true  }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Sno
    $ws.Cells.Item($r.Row, 2).Value = $r.Name
    $ws.Cells.Item($r.Row, 3).Value = $r.Value
    if ($r.Styled) {
        $ws.Cells.Item($r.Row, 2).NumberFormat = "@"
    }
}

# Update the selected range shown when the workbook is reopened.
[void]$ws.Range("B13:C13").Select()

Write-Output "applied"
